# Applies the "added trial to consolidated sheet / added mp's searsia model"
# edit: a new row (42) is appended to the "Allometric Models" sheet with the
# Searsia longispina canopy-area model, the print setup for that sheet is
# configured, and the active-sheet/selection bookkeeping is swapped from
# "Wet Dry Ratios" back to "Allometric Models".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Allometric Models"
$ws2 = $wb.Worksheets.Item(2)   # "Wet Dry Ratios"

# --- New row 42 on the Allometric Models sheet -----------------------------
# (order chosen so new shared-string entries land as MP, then the Log10
#  formula text, then CA)
$ws1.Cells.Item(42, 1).Value = "Searsia"
$ws1.Cells.Item(42, 2).Value = "longispina"
$ws1.Cells.Item(42, 3).Value = "MP"

$ws1.Cells.Item(42, 16).Font.Name = "Arial"
$ws1.Cells.Item(42, 16).Value = "Log10 y (C (kg) = 1.1012(Log10 canopy area (m2)) - 0.2938 "

$ws1.Cells.Item(42, 4).Value = "CA"
$ws1.Cells.Item(42, 5).Value = 24
$ws1.Cells.Item(42, 6).Value = 0.5077
$ws1.Cells.Item(42, 7).Value = 1.1012
$ws1.Cells.Item(42, 8).Value = -0.2938
$ws1.Cells.Item(42, 15).Value = "x"

# --- Print setup for the sheet now that it has data past the printable area
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# --- Active sheet / selection bookkeeping -----------------------------------
# Leave "Wet Dry Ratios" with its own last-used selection but not active...
$ws2.Activate()
$ws2.Range("I19").Select()
# ...then make "Allometric Models" the active/selected tab again.
$ws1.Activate()
$ws1.Range("H43").Select()
